$d = $word.ActiveDocument

# --- 1. Title paragraph: "PENGEMBANGAN SISTEM INFORMASI ..." -> "PERANCANGAN APLIKASI ..." ---
$pTitle = $d.Paragraphs(1).Range
if ($pTitle.Text -like "PENGEMBANGAN SISTEM INFORMASI*") {
    $pTitle.Find.Execute(
        "PENGEMBANGAN SISTEM INFORMASI", $false, $false, $false, $false, $false, $true, 1, $false,
        "PERANCANGAN APLIKASI", 2) | Out-Null
}

# --- 2. "Admin" role line -> append role description ---
$pAdmin = $d.Paragraphs(22).Range
if ($pAdmin.Text.TrimEnd([char]13) -eq "Admin") {
    $pAdmin.MoveEnd(1, -1) | Out-Null
    $pAdmin.InsertAfter(" (Petugas dari Dinas Pariwista yang mengelola aplikasi)")
}

# --- 3. "User" role line -> "Pengunjung / Wisatawan" ---
$pUser = $d.Paragraphs(25).Range
if ($pUser.Text.TrimEnd([char]13) -eq "User") {
    $pUser.MoveEnd(1, -1) | Out-Null
    $pUser.Text = "Pengunjung / Wisatawan"
}

# --- 4. "Pengelola" role line -> append " Tempat Wisata" ---
$pPengelola = $d.Paragraphs(30).Range
if ($pPengelola.Text.TrimEnd([char]13) -eq "Pengelola") {
    $pPengelola.MoveEnd(1, -1) | Out-Null
    $pPengelola.InsertAfter(" Tempat Wisata")
}

# --- 5. "Dinas Pariwisata " list item -> prepend "Pimpinan " ---
$pDinas1 = $d.Paragraphs(34).Range
if ($pDinas1.Text.TrimEnd([char]13) -eq "Dinas Pariwisata ") {
    $pDinas1.InsertBefore("Pimpinan ")
}

# --- 6. "Dinas pariwisata dapat menerima ..." list item -> prepend "Pimpinan " ---
$pDinas2 = $d.Paragraphs(35).Range
if ($pDinas2.Text.TrimEnd([char]13) -eq "Dinas pariwisata dapat menerima laporan dari setiap tempat wisata.") {
    $pDinas2.InsertBefore("Pimpinan ")
}
